# PowerBoard V3 Logbook - add new entry and tidy up the PCBnew comment
# that spilled across two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22/23 (2018-03-21, CB): the comment used to be split across two
# unmerged cells (B22 and B23). Combine it into a single sentence and
# merge the cells, matching the layout used by every other log entry.
$ws.Range("B23").Value = ""
$ws.Range("B22").Value = "Started PCBnew. changed footprints to my kicad's format. Layed out the general format for the pcb board. Added the IRF3205 transistor to the BOM. "
$ws.Range("B22:J23").Merge()

# --- New log entry: 2018-03-27 (CB) ---
# Use the existing 2018-03-24 entry (rows 24:25) as the formatting template
# for the new rows 26:27.
$ws.Range("A24:J25").Copy()
$ws.Range("A26:J27").PasteSpecial(-4122)

$ws.Range("A26").Value = 43187
$ws.Range("B26").Value = "Placed tracks for Motorcontrollers, Cells and the 9V IC. Check 9V to see if its right as I'm not sure"
$ws.Range("A27").Value = "CB"
$ws.Range("B27").Value = ""
$ws.Range("B26:J27").Merge()

$ws.Range("E29").Select()
